$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.893.15'
$ws.Range('E2').Value = '  +2.99%  '

$ws.Range('D3').Value = '2.304.18'
$ws.Range('E3').Value = '  +1.11%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.19'
$ws.Range('E5').Value = '  +2.47%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.82'
$ws.Range('E6').Value = '  +2.64%  '

$ws.Range('E7').Value = '  +1.29%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  +1.53%  '

$ws.Range('E10').Value = '  +4.92%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('E12').Value = '  +5.34%  '

$ws.Range('E13').Value = '  +0.95%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.978'
$ws.Range('E14').Value = '  +2.76%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.39'
$ws.Range('E15').Value = '  +1.21%  '

$ws.Range('D16').Value = '2.654.05'
$ws.Range('E16').Value = '  +1.33%  '

$ws.Range('D17').Value = '2.326.82'
$ws.Range('E17').Value = '  +2.20%  '

$ws.Range('D18').Value = '42.844.31'
$ws.Range('E18').Value = '  +2.66%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.57'
$ws.Range('E19').Value = '  +1.74%  '

$ws.Range('E20').Value = '  +1.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.28'
$ws.Range('E21').Value = '  +33.80%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.88'
$ws.Range('E22').Value = '  +1.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.62'
$ws.Range('E23').Value = '  +2.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '273.11'
$ws.Range('E24').Value = '  -1.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.24'
$ws.Range('E25').Value = '  +0.09%  '

$ws.Range('E26').Value = '  -0.55%  '

$ws.Range('E27').Value = '  +2.74%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.32'
$ws.Range('E28').Value = '  -0.43%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.73'
$ws.Range('E29').Value = '  -0.43%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.09'
$ws.Range('E30').Value = '  +10.46%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '165.85'
$ws.Range('E31').Value = '  +2.01%  '

$ws.Range('E32').Value = '  +6.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0891'
$ws.Range('E33').Value = '  +2.88%  '

$ws.Range('E34').Value = '  +1.13%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.115'
$ws.Range('E35').Value = '  +1.20%  '

$ws.Range('E36').Value = '  -11.22%  '

$ws.Range('E37').Value = '  +3.35%  '

$ws.Range('E38').Value = '  +3.53%  '

$ws.Range('E39').Value = '  +3.06%  '

$ws.Range('E40').Value = '  -4.29%  '

$ws.Range('E41').Value = '  +8.04%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.75'
$ws.Range('E42').Value = '  +3.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.23'
$ws.Range('E43').Value = '  +3.85%  '

$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.226'
$ws.Range('E44').Value = '  +1.61%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.49'
$ws.Range('E46').Value = '  +6.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '83.05'
$ws.Range('E47').Value = '  +10.92%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.96'
$ws.Range('E48').Value = '  -0.04%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.32'
$ws.Range('E49').Value = '  +1.69%  '

$ws.Range('E50').Value = '  -0.20%  '

$ws.Range('D51').Value = '1.593.28'
$ws.Range('E51').Value = '  +4.79%  '
